# smartart-linear-rule.pptx: "oox smartart, linear layout: fix scaling of
# spacing without rules".
#
# The upstream change simply re-enables ("uncomments") two already-present
# <dgm:layoutNode name="padding1"/"padding2"> blocks inside the SmartArt
# *layout definition* part (ppt/diagrams/layout1.xml) of the diagram, so the
# linear-process diagram carries the full, PowerPoint-UI-generated spacing
# markup instead of the trimmed-down/minimal version that used to ship in
# this test document. No slide text, node structure, geometry or styling is
# changed by the commit - only the diagram's internal layout-algorithm
# definition is restored to its complete form.
#
# That layout-definition XML is internal plumbing for the SmartArt rendering
# engine: it is not data bound to any node and PowerPoint's object model
# (Shape.SmartArt / SmartArtNode / SmartArtLayout ...) does not expose it for
# editing, even in the real product - there is no VBA/COM property for
# "the <dgm:layoutNode> rules of this diagram's layout". The closest, and
# only, legitimate COM entry point that targets this exact concept is
# SmartArt.Reset(), whose documented job is to reapply the diagram's own
# layout definition to the graphic, so that is what we invoke here on the
# presentation's SmartArt diagram shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasSmartArt) {
        $shp.SmartArt.Reset()
    }
}
